# Auto-generated: applies exact cell changes from the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.316.24"
$ws.Range("E2").Value = "  -1.44%  "
# Row 3
$ws.Range("D3").Value = "1.828.71"
$ws.Range("E3").Value = "  -1.24%  "
# Row 4
$ws.Range("E4").Value = "  -0.77%  "
# Row 5
$ws.Range("D5").Value = "314.55"
$ws.Range("E5").Value = "  -1.59%  "
# Row 6
$ws.Range("E6").Value = "  -0.57%  "
# Row 7
$ws.Range("D7").Value = "0.4237"
$ws.Range("E7").Value = "  -2.04%  "
# Row 8
$ws.Range("D8").Value = "0.3693"
$ws.Range("E8").Value = "  -1.91%  "
# Row 9
$ws.Range("D9").Value = "0.07266"
$ws.Range("E9").Value = "  -1.70%  "
# Row 10
$ws.Range("D10").Value = "0.8671"
$ws.Range("E10").Value = "  -2.07%  "
# Row 11
$ws.Range("E11").Value = "  -2.65%  "
# Row 12
$ws.Range("D12").Value = "1.824.46"
$ws.Range("E12").Value = "  -1.72%  "
# Row 13
$ws.Range("D13").Value = "6.742"
$ws.Range("E13").Value = "  -0.29%  "
# Row 14
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "0.07095"
$ws.Range("E14").Value = "  -0.25%  "
# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "5.324"
$ws.Range("E15").Value = "  -2.90%  "
# Row 16
$ws.Range("D16").Value = "89.48"
$ws.Range("E16").Value = "  +1.33%  "
# Row 17
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  -0.86%  "
# Row 18
$ws.Range("D18").Value = "0.000008878"
$ws.Range("E18").Value = "  -1.81%  "
# Row 19
$ws.Range("E19").Value = "  -0.64%  "
# Row 20
$ws.Range("B20").Value = "BitDAO"
$ws.Range("C20").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D20").Value = "0.5063"
$ws.Range("E20").Value = "  +2.88%  "
# Row 21
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "15.13"
$ws.Range("E21").Value = "  -2.55%  "
# Row 22
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "27.361.81"
$ws.Range("E22").Value = "  -1.29%  "
# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.141"
$ws.Range("E23").Value = "  -2.60%  "
# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  -2.63%  "
# Row 25
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.051.73"
$ws.Range("E25").Value = "  -2.55%  "
# Row 26
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "1.996"
$ws.Range("E26").Value = "  -1.79%  "
# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "152.98"
$ws.Range("E27").Value = "  -1.86%  "
# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.181"
$ws.Range("E28").Value = "  +2.08%  "
# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "18.40"
$ws.Range("E29").Value = "  -1.33%  "
# Row 30
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "5.251"
$ws.Range("E30").Value = "  -3.36%  "
# Row 31
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "116.54"
$ws.Range("E31").Value = "  -3.72%  "
# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.08879"
$ws.Range("E32").Value = "  -0.99%  "
# Row 33
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "1.208"
$ws.Range("E33").Value = "  -2.73%  "
# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.7604"
$ws.Range("E34").Value = "  -3.10%  "
# Row 35
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "4.475"
$ws.Range("E35").Value = "  -2.38%  "
# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.833"
$ws.Range("E36").Value = "  -3.05%  "
# Row 37
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "1.005"
$ws.Range("E37").Value = "  -0.71%  "
# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.122"
$ws.Range("E38").Value = "  -2.12%  "
# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01981"
$ws.Range("E39").Value = "  +0.51%  "
# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.05278"
$ws.Range("E40").Value = "  -1.18%  "
# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "7.293"
$ws.Range("E41").Value = "  +2.24%  "
# Row 42
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "2.883"
$ws.Range("E42").Value = "  +0.53%  "
# Row 43
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.1702"
$ws.Range("E43").Value = "  +0.81%  "
# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.5080"
$ws.Range("E44").Value = "  -2.21%  "
# Row 45
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "8.710"
$ws.Range("E45").Value = "  -4.12%  "
# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "10.71"
$ws.Range("E46").Value = "  -0.32%  "
# Row 47
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "107.93"
$ws.Range("E47").Value = "  -2.64%  "
# Row 48
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "0.4773"
$ws.Range("E48").Value = "  +0.72%  "
# Row 49
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "1.005"
$ws.Range("E49").Value = "  -0.58%  "
# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.06390"
$ws.Range("E50").Value = "  -1.93%  "
# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.672"
$ws.Range("E51").Value = "  -2.94%  "
